$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell replacements
$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"
$t.Cell(4,1).Range.Text = "473"

$t.Cell(6,1).Range.Text = "1.00193"
$t.Cell(7,1).Range.Text = "0.16183"
$t.Cell(8,1).Range.Text = "0.05918"

$t.Cell(11,1).Range.Text = "0.96647"
$t.Cell(12,1).Range.Text = "24.67237"

# Collapse the multi-run, tab-separated cells into single values
$t.Cell(44,1).Range.Text = "87.66"
$t.Cell(45,1).Range.Text = "24.67"
$t.Cell(46,1).Range.Text = "199"
